$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Issue List")

$xlPasteFormats = [Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteFormats

# Resolved some issues opened on 2014/11/10:
# set 状态 (status, column D) to 已解决 (Resolved) and fill in
# 解决日期 (resolution date, column F) for rows 10-16.
$rows = @(
    @{ Row = 10; Date = 41954 },
    @{ Row = 11; Date = 41955 },
    @{ Row = 12; Date = 41954 },
    @{ Row = 13; Date = 41954 },
    @{ Row = 14; Date = 41955 },
    @{ Row = 15; Date = 41954 },
    @{ Row = 16; Date = 41955 }
)

foreach ($item in $rows) {
    $r = $item.Row

    $ws.Cells.Item($r, 4).Value = "已解决"

    # Column F has no date format yet on these rows; copy the existing
    # date formatting from column E (same row) before writing the value.
    $ws.Cells.Item($r, 5).Copy()
    $ws.Cells.Item($r, 6).PasteSpecial($xlPasteFormats)
    $ws.Cells.Item($r, 6).Value = $item.Date
}

$excel.CutCopyMode = 0

# Leave the view where the author ended up editing.
$ws.Application.ActiveWindow.ScrollRow = 7
$ws.Range("F16").Select() | Out-Null
